$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(8).Insert()
$ws.Cells.Item(1, 8).Value = "diploma"

$ws.Range("H1").Select() | Out-Null

